$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking row: correct answers count
$ws.Range("B11").Value = 5

# Update total row: corrected total marks and corr/total string
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
